$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13.67106065199291
$ws.Range("E2").Value = 8.810827213099513
$ws.Range("H2").Value = 13.69137680584594
$ws.Range("I2").Value = 8.854824500149302
